# Helper approach: for D-column price values that look numeric (e.g. "64.49"),
# Excel would auto-convert them to a Number type and attach a new cell style.
# To preserve them as plain text cells (matching the source data format),
# we temporarily force a Text number format, assign the value, then reset
# the cell style back to Normal/default so no stray style survives.

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "26.711.42"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "1.638.87"
$ws.Range("E3").Value = "  -0.60%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "217.91"
$ws.Range("E5").Value = "  +0.51%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.83%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.21%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.72%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.87%  "

# Row 10 - Solana
Set-TextValue $ws "D10" "19.07"
$ws.Range("E10").Value = "  -0.63%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.00%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D12" "1.866.68"
$ws.Range("E12").Value = "  -0.71%  "

# Row 13 - WrappedEther
Set-TextValue $ws "D13" "1.640.81"
$ws.Range("E13").Value = "  -0.43%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.36%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -1.61%  "

# Row 16 - Litecoin
Set-TextValue $ws "D16" "64.49"
$ws.Range("E16").Value = "  -1.49%  "

# Row 17 - WrappedBTC
Set-TextValue $ws "D17" "26.697.10"

# Row 18 - ShibaInu
Set-TextValue $ws "D18" "0.0₃0727"
$ws.Range("E18").Value = "  -2.38%  "

# Row 19 - BitcoinCash
Set-TextValue $ws "D19" "211.23"
$ws.Range("E19").Value = "  -3.29%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.16%  "

# Row 22 - Chainlink
Set-TextValue $ws "D22" "6.18"

# Row 23 - Toncoin
Set-TextValue $ws "D23" "2.30"
$ws.Range("E23").Value = "  -4.49%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  -2.61%  "

# Row 25 - Monero
Set-TextValue $ws "D25" "146.73"
$ws.Range("E25").Value = "  +0.19%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.08%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -1.89%  "

# Row 28 - Cosmos
Set-TextValue $ws "D28" "7.08"
$ws.Range("E28").Value = "  -0.78%  "

# Row 29 - EthereumClassic
Set-TextValue $ws "D29" "15.56"
$ws.Range("E29").Value = "  -1.08%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -3.04%  "

# Row 31 - PancakeSwap
Set-TextValue $ws "D31" "1.18"
$ws.Range("E31").Value = "  +0.23%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.21%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.92%  "

# Row 34 - Maker
Set-TextValue $ws "D34" "1.265.89"
$ws.Range("E34").Value = "  -1.50%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -1.06%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.73%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -2.03%  "

# Row 38 - ImmutableX
Set-TextValue $ws "D38" "0.528"
$ws.Range("E38").Value = "  -1.83%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  -2.91%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  -0.26%  "

# Row 41 - TrustWalletToken
Set-TextValue $ws "D41" "0.802"

# Row 42 - MXToken
$ws.Range("E42").Value = "  -3.14%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -3.68%  "

# Row 44 - RocketPoolETH
Set-TextValue $ws "D44" "1.776.60"
$ws.Range("E44").Value = "  -0.74%  "

# Row 45 - Quant
Set-TextValue $ws "D45" "91.34"
$ws.Range("E45").Value = "  -0.80%  "

# Row 46 - Aave
$ws.Range("E46").Value = "  +0.66%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -2.03%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +0.47%  "

# Row 49 - EnergySwap
Set-TextValue $ws "D49" "7.51"
$ws.Range("E49").Value = "  -3.48%  "

# Row 50 - Algorand
Set-TextValue $ws "D50" "0.0959"
$ws.Range("E50").Value = "  -1.14%  "

# Row 51 - USDD -> Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D51" "0.407"
$ws.Range("E51").Value = "  -0.41%  "
